# V10.0 - zmienna rozdzielczosc (27 nowych wierszy wynikow) + tlumaczenie
# naglowkow i nazw wojewodztw na polskie (poza tytulem arkusza)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 27 fresh rows right below the header, pushing the old ---
#        24 data rows down from 2..26 to 29..53.
$ws.Range("A2:A28").EntireRow.Insert()

# The inserted rows inherit the header's bold/bordered style on B:E;
# strip that back out so B2:E28 start out unstyled, same as the other
# data rows in the sheet.
$ws.Range("B2:E28").ClearFormats()

# Give column A of the new rows the same style used by the rest of
# column A (bold/bordered, centered) by copying the format from a
# cell that still has it (now at row 29) down onto A2:A28.
$ws.Range("A29").Copy()
$ws.Range("A2:A28").PasteSpecial(-4122)

# --- 2. Translate the header row into Polish (title bar/tab name is ---
#        left alone, only the column headers change).
$ws.Cells.Item(1, 2).Value = "Nazwa"
$ws.Cells.Item(1, 3).Value = "Poziom"
$ws.Cells.Item(1, 4).Value = "Województwo"
$ws.Cells.Item(1, 5).Value = "Punkty"

# --- 3. Fill in the 27 new result rows (B=name, C=level, D=voivodeship, E=points).
$newRows = @(
    @("maks", "Extreme", "Podlaskie", 9),
    @("maks", "Extreme", "Podlaskie", 12),
    @("kk", "Extreme", "Kujawsko-Pomorskie", 12),
    @("h", "Extreme", "Podlaskie", 9),
    @("hh", "Medium", "Podlaskie", 3),
    @("ii", "Medium", "Małopolskie", 5),
    @("aa", "Extreme", "Podlaskie", 15),
    @("aaa", "Extreme", "Kujawsko-Pomorskie", 21),
    @("a", "Extreme", "Podlaskie", 9),
    @("aaa", "Extreme", "Kujawsko-Pomorskie", 18),
    @("jjj", "Extreme", "Opolskie", 15),
    @(";;", "Extreme", "Dolnośląskie", 21),
    @("aaa", "Extreme", "Łódzkie", 24),
    @("aa", "Extreme", "Dolnośląskie", 27),
    @("", "Extreme", "Podlaskie", 15),
    @("", "Extreme", "Łódzkie", 27),
    @("lala", "Extreme", "Kujawsko-Pomorskie", 18),
    @("kuba", "Extreme", "Podlaskie", 12),
    @("Maks", "Extreme", "Opolskie", 9),
    @("Maks", "Extreme", "Opolskie", 6),
    @("Maks", "Extreme", "Opolskie", 9),
    @("Maks", "Extreme", "Opolskie", 15),
    @("Maks", "Extreme", "Opolskie", 12),
    @("Maks", "Extreme", "Opolskie", 3),
    @("Maks", "Extreme", "Opolskie", 6),
    @("Maks", "Extreme", "Łódzkie", 15),
    @("Karolcio", "Extreme", "Wielkopolskie", 42)
)

$rowNum = 2
foreach ($row in $newRows) {
    # Skip blank "name" entries: they must stay as an empty cell (same
    # as every other blank-name row already in the sheet), and writing
    # an explicit "" would not round-trip the same way as never having
    # written anything there.
    if ($row[0] -ne "") {
        $ws.Cells.Item($rowNum, 2).Value = $row[0]
    }
    $ws.Cells.Item($rowNum, 3).Value = $row[1]
    $ws.Cells.Item($rowNum, 4).Value = $row[2]
    $ws.Cells.Item($rowNum, 5).Value = $row[3]
    $rowNum = $rowNum + 1
}

# --- 4. Column A is just a running 0-based index for every data row; ---
#        after the insert it is wrong for the inserted block (blank)
#        and needs no change for the shifted rows (Excel doesn't
#        renumber them automatically), so recompute A for every row.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

Write-Host "Done: rows now " $ws.UsedRange.Address()
